# Implements "implemented reading scenario data from csv":
#   - ScenarioTexts gets new columns: ScenarioDay, ScenarioScope, ScenarioItemCount,
#     ScenarioItem1..4 and Champion1..4, replacing the old single "ScenarioItems"
#     ("Yes, No") column.
#   - ScenarioAttributes header row is untouched in content (ScenarioId, ScenarioDay,
#     ScenarioScope, ScenarioItems); only its selection changes.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("ScenarioTexts")
$ws2 = $wb.Worksheets.Item("ScenarioAttributes")

# ---------------------------------------------------------------------------
# Sheet 1: ScenarioTexts
# ---------------------------------------------------------------------------

# Header row
$ws1.Cells.Item(1,1).Value  = "ScenarioId"
$ws1.Cells.Item(1,2).Value  = "ScenarioName"
$ws1.Cells.Item(1,3).Value  = "ScenarioDescription"
$ws1.Cells.Item(1,4).Value  = "ScenarioDay"
$ws1.Cells.Item(1,5).Value  = "ScenarioScope"
$ws1.Cells.Item(1,6).Value  = "ScenarioItemCount"
$ws1.Cells.Item(1,7).Value  = "ScenarioItem1"
$ws1.Cells.Item(1,8).Value  = "ScenarioItem2"
$ws1.Cells.Item(1,9).Value  = "ScenarioItem3"
$ws1.Cells.Item(1,10).Value = "ScenarioItem4"
$ws1.Cells.Item(1,11).Value = "Champion1"
$ws1.Cells.Item(1,12).Value = "Champion2"
$ws1.Cells.Item(1,13).Value = "Champion3"
$ws1.Cells.Item(1,14).Value = "Champion4"

# Row 2: Prologue
$ws1.Cells.Item(2,1).Value  = 1
$ws1.Cells.Item(2,2).Value  = "Prologue"
$ws1.Cells.Item(2,3).Value  = "Marry partner"
$ws1.Cells.Item(2,4).Value  = 1
$ws1.Cells.Item(2,5).Value  = 1
$ws1.Cells.Item(2,6).Value  = 2
$ws1.Cells.Item(2,7).Value  = "Yes"
$ws1.Cells.Item(2,8).Value  = "No"
$ws1.Cells.Item(2,11).Value = "Ragnar"
$ws1.Cells.Item(2,12).Value = "Lagertha"

# Row 3: Epilogue
$ws1.Cells.Item(3,1).Value  = 2
$ws1.Cells.Item(3,2).Value  = "Epilogue"
$ws1.Cells.Item(3,3).Value  = "Arrange funeral"
$ws1.Cells.Item(3,4).Value  = 2
$ws1.Cells.Item(3,5).Value  = 1
$ws1.Cells.Item(3,6).Value  = 4
$ws1.Cells.Item(3,7).Value  = "Yes"
$ws1.Cells.Item(3,8).Value  = "No"
$ws1.Cells.Item(3,9).Value  = "Yes, but…"
$ws1.Cells.Item(3,10).Value = "No, but…"
$ws1.Cells.Item(3,11).Value = "Ragnar"
$ws1.Cells.Item(3,12).Value = "Widow"
$ws1.Cells.Item(3,13).Value = "Deceased"

# Column widths for the newly-populated columns (bestfit-style, matches the
# behaviour already applied to columns B/C by the original author)
$ws1.Columns.Item(5).AutoFit()  | Out-Null
$ws1.Columns.Item(6).AutoFit()  | Out-Null
$ws1.Columns.Item(7).AutoFit()  | Out-Null
$ws1.Columns.Item(8).AutoFit()  | Out-Null
$ws1.Columns.Item(9).AutoFit()  | Out-Null
$ws1.Columns.Item(10).AutoFit() | Out-Null

# ---------------------------------------------------------------------------
# Sheet 2: ScenarioAttributes (content unchanged, just re-asserted)
# ---------------------------------------------------------------------------

$ws2.Cells.Item(1,1).Value = "ScenarioId"
$ws2.Cells.Item(1,2).Value = "ScenarioDay"
$ws2.Cells.Item(1,3).Value = "ScenarioScope"
$ws2.Cells.Item(1,4).Value = "ScenarioItems"

# ---------------------------------------------------------------------------
# Selections: ScenarioAttributes first so ScenarioTexts ends up the active tab
# ---------------------------------------------------------------------------

$ws2.Range("A1:D1").Select() | Out-Null
$ws1.Activate() | Out-Null
$ws1.Range("N2").Select() | Out-Null
